# Generate Report for Handoff
# Updates the handoff/xliff-generation timestamps for file
# 33d08fc1-2d12-460d-bdad-c6841c280200 across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Row 6 on every sheet corresponds to 33d08fc1-2d12-460d-bdad-c6841c280200

# Overview sheet: "Latest HO Xliff Generate Date" column G
$wsOverview.Range("G6").Value = "2016-08-18 12:42:08"

# zh-cn sheet: "Latest Handoff Datetime" column H
$wsZhCn.Range("H6").Value = "2016-08-18 12:41:57"

# de-de sheet: "Latest Handoff Datetime" column H
$wsDeDe.Range("H6").Value = "2016-08-18 12:42:08"
